$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 53
$ws.Range("H53").Value = 22223428
$ws.Range("I53").Value = 66667700
$ws.Range("J53").Value = 1291.8
$ws.Range("K53").Value = 66667700
$ws.Range("L53").Value = 1291.8
$ws.Range("M53").Value = -66667063
$ws.Range("N53").Value = -2565.8

# ALC row 86
$ws.Range("H86").Value = 3764019.8
$ws.Range("I86").Value = 3783.6667
$ws.Range("J86").Value = 6584197
$ws.Range("K86").Value = 3783.6667
$ws.Range("L86").Value = 6584197
$ws.Range("M86").Value = -2660.6667
$ws.Range("N86").Value = -6586443

# ALC row 89
$ws.Range("H89").Value = 3764019.8
$ws.Range("I89").Value = 3783.6667
$ws.Range("J89").Value = 6584197
$ws.Range("K89").Value = 18918.3335
$ws.Range("L89").Value = 32920985
$ws.Range("M89").Value = -13302.3335
$ws.Range("N89").Value = -32932217

# ALC row 128
$ws.Range("H128").Value = 89393.914
$ws.Range("J128").Value = 89393.914
$ws.Range("L128").Value = 89393.914
$ws.Range("N128").Value = -99353.914

# ALC row 132
$ws.Range("H132").Value = 2142.5425
$ws.Range("I132").Value = 1970.3112
$ws.Range("K132").Value = 5910.9336
$ws.Range("M132").Value = -3380.9336

# ALC row 135
$ws.Range("H135").Value = 513915.7
$ws.Range("I135").Value = 589078.9
$ws.Range("J135").Value = 2806.2
$ws.Range("K135").Value = 5301710.100000001
$ws.Range("L135").Value = 25255.8
$ws.Range("M135").Value = -5299175.100000001
$ws.Range("N135").Value = -30325.8

# ALC row 137
$ws.Range("H137").Value = 367505.3
$ws.Range("I137").Value = 235770.52
$ws.Range("K137").Value = 707311.5599999999
$ws.Range("M137").Value = -704761.5599999999

# ALC row 138
$ws.Range("H138").Value = 4832.808
$ws.Range("I138").Value = 2401.8667
$ws.Range("J138").Value = 6528.814
$ws.Range("K138").Value = 7205.6001
$ws.Range("L138").Value = 19586.442
$ws.Range("M138").Value = -2065.6001
$ws.Range("N138").Value = -29866.442

# ALC row 141
$ws.Range("H141").Value = 2566.691
$ws.Range("I141").Value = 1167.561
$ws.Range("K141").Value = 3502.683
$ws.Range("M141").Value = 1677.317

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 2748.6626
$ws.Range("I32").Value = 1806.3206
$ws.Range("K32").Value = 1806.3206
$ws.Range("M32").Value = -1519.3206

# ARM row 61
$ws.Range("H61").Value = 3213.182
$ws.Range("I61").Value = 2766.111
$ws.Range("K61").Value = 2766.111
$ws.Range("M61").Value = -2554.111

# ARM row 102
$ws.Range("H102").Value = 1323
$ws.Range("I102").Value = 1323
$ws.Range("K102").Value = 1323
$ws.Range("M102").Value = 299

# ARM row 110
$ws.Range("H110").Value = 1172.8
$ws.Range("I110").Value = 1068.7307
$ws.Range("K110").Value = 1068.7307
$ws.Range("M110").Value = 976.2692999999999

# ARM row 132
$ws.Range("H132").Value = 3606.4285
$ws.Range("I132").Value = 2390.4783
$ws.Range("J132").Value = 9199.799999999999
$ws.Range("K132").Value = 7171.4349
$ws.Range("L132").Value = 27599.4
$ws.Range("M132").Value = -4641.4349
$ws.Range("N132").Value = -32659.4

# ARM row 136
$ws.Range("H136").Value = 3213.182
$ws.Range("I136").Value = 2766.111
$ws.Range("K136").Value = 8298.332999999999
$ws.Range("M136").Value = -5748.332999999999

$ws = $wb.Worksheets.Item("BSM")
# BSM row 5
$ws.Range("H5").Value = 4003.3333
$ws.Range("I5").Value = 4003.3333
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4003.3333
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null

# BSM row 75
$ws.Range("H75").Value = 13250
$ws.Range("I75").Value = 13250
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 13250
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = $null

# BSM row 78
$ws.Range("H78").Value = 13250
$ws.Range("I78").Value = 13250
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 39750
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = $null

# BSM row 81
$ws.Range("H81").Value = 42913.332
$ws.Range("J81").Value = 42913.332
$ws.Range("L81").Value = 42913.332
$ws.Range("N81").Value = -45035.332

# BSM row 84
$ws.Range("H84").Value = 42913.332
$ws.Range("J84").Value = 42913.332
$ws.Range("L84").Value = 128739.996
$ws.Range("N84").Value = -139347.996

# BSM row 105
$ws.Range("H105").Value = 3130.5557
$ws.Range("I105").Value = 2739.5715
$ws.Range("K105").Value = 2739.5715
$ws.Range("M105").Value = -992.5715

# BSM row 107
$ws.Range("H107").Value = 419548.4
$ws.Range("I107").Value = 2629.2104
$ws.Range("J107").Value = 2003841.4
$ws.Range("K107").Value = 2629.2104
$ws.Range("L107").Value = 2003841.4
$ws.Range("M107").Value = -709.2103999999999
$ws.Range("N107").Value = -2007681.4

# BSM row 134
$ws.Range("H134").Value = 26424.455
$ws.Range("I134").Value = 3566.6
$ws.Range("J134").Value = 115316.11
$ws.Range("K134").Value = 10699.8
$ws.Range("L134").Value = 345948.33
$ws.Range("M134").Value = -8164.799999999999
$ws.Range("N134").Value = -351018.33

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 61090.055
$ws.Range("I31").Value = 1690.1428
$ws.Range("J31").Value = 98890
$ws.Range("K31").Value = 1690.1428
$ws.Range("L31").Value = 98890
$ws.Range("M31").Value = -1395.1428
$ws.Range("N31").Value = -99480

# CRP row 34
$ws.Range("H34").Value = 61090.055
$ws.Range("I34").Value = 1690.1428
$ws.Range("J34").Value = 98890
$ws.Range("K34").Value = 1690.1428
$ws.Range("L34").Value = 98890
$ws.Range("M34").Value = -1488.1428
$ws.Range("N34").Value = -99294

# CRP row 132
$ws.Range("H132").Value = 1876.7108
$ws.Range("I132").Value = 1578.7778
$ws.Range("K132").Value = 4736.3334
$ws.Range("M132").Value = -2206.3334

$ws = $wb.Worksheets.Item("CUL")
# CUL row 23
$ws.Range("H23").Value = 833.3333
$ws.Range("I23").Value = 500
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 3000

# CUL row 106
$ws.Range("H106").Value = 16597
$ws.Range("J106").Value = 13194
$ws.Range("L106").Value = 39582
$ws.Range("N106").Value = -41474

# CUL row 113
$ws.Range("H113").Value = 2059041.6
$ws.Range("J113").Value = 1576
$ws.Range("L113").Value = 4728
$ws.Range("N113").Value = -9068

# CUL row 132
$ws.Range("H132").Value = 848322.6
$ws.Range("J132").Value = 911426.8
$ws.Range("L132").Value = 8202841.2
$ws.Range("N132").Value = -8207901.2

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 1824009.1
$ws.Range("I80").Value = 1433842.6
$ws.Range("K80").Value = 1433842.6
$ws.Range("M80").Value = -1432844.6

# GSM row 83
$ws.Range("H83").Value = 1824009.1
$ws.Range("I83").Value = 1433842.6
$ws.Range("K83").Value = 7169213
$ws.Range("M83").Value = -7164221

# GSM row 94
$ws.Range("H94").Value = 42854
$ws.Range("I94").Value = 41000
$ws.Range("J94").Value = 43163
$ws.Range("K94").Value = 41000
$ws.Range("L94").Value = 43163
$ws.Range("N94").Value = -44515

# GSM row 120
$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000

# GSM row 122
$ws.Range("H122").Value = 4062.8333
$ws.Range("I122").Value = 3935.4
$ws.Range("K122").Value = 11806.2
$ws.Range("M122").Value = -9356.200000000001

# GSM row 132
$ws.Range("H132").Value = 585178.6
$ws.Range("I132").Value = 671143
$ws.Range("J132").Value = 262812.25
$ws.Range("K132").Value = 2013429
$ws.Range("L132").Value = 788436.75
$ws.Range("M132").Value = -2010899
$ws.Range("N132").Value = -793496.75

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Range("H40").Value = 1746.1714
$ws.Range("I40").Value = 1627.2903
$ws.Range("J40").Value = 2667.5
$ws.Range("K40").Value = 1627.2903
$ws.Range("L40").Value = 2667.5
$ws.Range("M40").Value = -1491.2903
$ws.Range("N40").Value = -2939.5

# LTW row 46
$ws.Range("H46").Value = 3494.2173
$ws.Range("I46").Value = 3312.1667
$ws.Range("J46").Value = 3692.818
$ws.Range("K46").Value = 3312.1667
$ws.Range("L46").Value = 3692.818
$ws.Range("M46").Value = -3124.1667
$ws.Range("N46").Value = -4068.818

# LTW row 82
$ws.Range("H82").Value = 2438.5
$ws.Range("I82").Value = 2438.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2438.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = $null

# LTW row 85
$ws.Range("H85").Value = 2438.5
$ws.Range("I85").Value = 2438.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2438.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = $null

$ws = $wb.Worksheets.Item("WVR")
# WVR row 8
$ws.Range("H8").Value = 19333.334
$ws.Range("I8").Value = 19333.334
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 19333.334
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = $null

# WVR row 81
$ws.Range("H81").Value = 2658.8
$ws.Range("I81").Value = 2823.75
$ws.Range("J81").Value = 1999
$ws.Range("K81").Value = 5647.5
$ws.Range("L81").Value = 3998
$ws.Range("M81").Value = -4586.5
$ws.Range("N81").Value = -6120

# WVR row 84
$ws.Range("H84").Value = 2658.8
$ws.Range("I84").Value = 2823.75
$ws.Range("J84").Value = 1999
$ws.Range("K84").Value = 28237.5
$ws.Range("L84").Value = 19990
$ws.Range("M84").Value = -22933.5
$ws.Range("N84").Value = -30598

# WVR row 107
$ws.Range("H107").Value = 1246.3684
$ws.Range("I107").Value = 1457.6666
$ws.Range("J107").Value = 454
$ws.Range("K107").Value = 4372.9998
$ws.Range("L107").Value = 1362
$ws.Range("M107").Value = -2452.9998
$ws.Range("N107").Value = -5202

# WVR row 132
$ws.Range("H132").Value = 22807.428
$ws.Range("I132").Value = 1148.683
$ws.Range("K132").Value = 3446.049
$ws.Range("M132").Value = -916.049

# WVR row 136
$ws.Range("H136").Value = 279516.8
$ws.Range("I136").Value = 326376.53
$ws.Range("J136").Value = 164169.84
$ws.Range("K136").Value = 979129.5900000001
$ws.Range("L136").Value = 492509.52
$ws.Range("M136").Value = -976579.5900000001
$ws.Range("N136").Value = -497609.52
